$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

$ws.Range("A11").Value = "paginationInsightName"
$ws.Range("B11").Value = "Automation_Insight"
$ws.Range("A12").Value = "Automation_Pivot_PaginationDashboard_DefaultPageSize"
$ws.Range("B12").Value = "Automation_Pivot_PaginationDashboard_DefaultPageSize"
$ws.Range("A13").Value = "Automation_Pivot_PaginationDashboard_PageSize20"
$ws.Range("B13").Value = "Automation_Pivot_PaginationDashboard_PageSize20"

$ws.Columns.Item(1).ColumnWidth = 48.83
